$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Utopía) updates
$ws.Range("T2").Value = 100
$ws.Range("W2").Value = 1.45508729917437
$ws.Range("AA2").Value = 1.56697332086377
$ws.Range("AB2").Value = -22.7103420555846
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 0
$ws.Range("AO2").Value = 100
$ws.Range("AU2").Value = 0
$ws.Range("AV2").Value = 0.770517376477767
$ws.Range("AW2").Value = 0
$ws.Range("AZ2").Value = 152.521817513277
$ws.Range("BC2").Value = 9.55150557584923
$ws.Range("BD2").Value = -83.0411921843169
$ws.Range("BE2").Value = 2724.90891581896

# Row 3 (Distopía) updates
$ws.Range("T3").Value = 0
$ws.Range("W3").Value = -9.44128729917438
$ws.Range("AA3").Value = -2.95067332086377
$ws.Range("AB3").Value = -27.8896579444154
$ws.Range("AI3").Value = 100
$ws.Range("AJ3").Value = 100
$ws.Range("AO3").Value = 0
$ws.Range("AU3").Value = 100
$ws.Range("AV3").Value = 0
$ws.Range("AW3").Value = 100
$ws.Range("AZ3").Value = 12.078182486723
$ws.Range("BC3").Value = -52.0332055758492
$ws.Range("BD3").Value = -134.435607815683
$ws.Range("BE3").Value = -373.908915818964
